# Update the answer cells in the single table on the page.
#
# Each cell is addressed explicitly by (row, col) and its text is
# replaced via a direct Range.Text assignment (scoped to the cell's own
# Start/End, minus the trailing end-of-cell mark) rather than
# Range.Find.Execute. Several source values repeat verbatim in more
# than one cell (e.g. "76÷3=25, 1" appears twice, mapping to two
# different targets), and this engine's Find scopes to the whole
# document rather than to the calling Range, so a plain Find/Replace
# would clobber the wrong cell. Direct Range.Text assignment is scoped
# correctly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellAnswer($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # $r.End sits just past the cell's end-of-cell mark (and paragraph
    # mark); trim 1 char so we only overwrite the visible text, leaving
    # the cell-structure markers intact.
    $inner = $d.Range($r.Start, $r.End - 1)
    $inner.Text = $newText
}

Set-CellAnswer $t 1 1 "86÷9=9, 5"
Set-CellAnswer $t 1 2 "87÷2=43, 1"
Set-CellAnswer $t 1 3 "45÷5=9, 0"
Set-CellAnswer $t 1 4 "49÷7=7, 0"
Set-CellAnswer $t 1 5 "97÷8=12, 1"

Set-CellAnswer $t 5 1 "33÷5=6, 3"
Set-CellAnswer $t 5 2 "87÷7=12, 3"
Set-CellAnswer $t 5 3 "92÷7=13, 1"
Set-CellAnswer $t 5 4 "52÷9=5, 7"
Set-CellAnswer $t 5 5 "64÷6=10, 4"

Set-CellAnswer $t 9 1 "48÷5=9, 3"
Set-CellAnswer $t 9 2 "99÷9=11, 0"
Set-CellAnswer $t 9 3 "18÷4=4, 2"
Set-CellAnswer $t 9 4 "30÷4=7, 2"
Set-CellAnswer $t 9 5 "26÷3=8, 2"

Set-CellAnswer $t 13 1 "40÷9=4, 4"
Set-CellAnswer $t 13 2 "25÷4=6, 1"
Set-CellAnswer $t 13 3 "10÷2=5, 0"
Set-CellAnswer $t 13 4 "62÷2=31, 0"
Set-CellAnswer $t 13 5 "82÷5=16, 2"

Set-CellAnswer $t 17 1 "59÷7=8, 3"
Set-CellAnswer $t 17 2 "98÷3=32, 2"
Set-CellAnswer $t 17 3 "80÷4=20, 0"
Set-CellAnswer $t 17 4 "69÷4=17, 1"
Set-CellAnswer $t 17 5 "75÷4=18, 3"

Write-Host "Done updating answer table."
